# RPA datasets push 2024-03-28
# Insert a new record row (하나32호스팩) at row 17, shifting the 한국/삼현
# and 한화/이에이트 rows down to 18/19 (their values stay identical, the
# Insert() shift handles that automatically). Also fix up rows 3/4/5/8
# which are reordered within the DB / NH groups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold date-like text (청약일/납입일/상장일 = B/F/G). These must
# stay plain text (shared-string) cells, not get auto-converted to Excel
# date serials, so force a Text format before writing into them and restore
# General afterwards. (Rows 1-18 cover every row we touch before the insert;
# the new row 17 is formatted again individually after the insert.)
$ws.Range("B2:B18").NumberFormat = "@"
$ws.Range("F2:F18").NumberFormat = "@"
$ws.Range("G2:G18").NumberFormat = "@"

function Set-Row($r, $vals) {
    $ws.Cells.Item($r, 1).Value  = $vals[0]
    $ws.Cells.Item($r, 2).Value  = $vals[1]
    $ws.Cells.Item($r, 3).Value  = $vals[2]
    $ws.Cells.Item($r, 4).Value  = $vals[3]
    $ws.Cells.Item($r, 5).Value  = $vals[4]
    $ws.Cells.Item($r, 6).Value  = $vals[5]
    $ws.Cells.Item($r, 7).Value  = $vals[6]
    $ws.Cells.Item($r, 8).Value  = $vals[7]
    $ws.Cells.Item($r, 9).Value  = $vals[8]
    $ws.Cells.Item($r, 10).Value = $vals[9]
    $ws.Cells.Item($r, 11).Value = $vals[10]
    $ws.Cells.Item($r, 12).Value = $vals[11]
}

# --- Fix rows 3, 4, 5, 8 (data got reshuffled within the DB/NH groups) ---
Set-Row 3 @("DB", "2024-01-25", "스튜디오삼익", "DB", "DB", "2024-01-30", "2024-02-06", 15300, 850000, 18000, 0, 100)
Set-Row 4 @("DB", "2024-02-26", "케이엔알시스템", "DB, NH", "DB, NH", "2024-02-29", "2024-03-07", 14202, 2104000, 13500, 0, 50)
Set-Row 5 @("NH", "2024-03-14", "엔젤로보틱스", "NH", "NH", "2024-03-19", "2024-03-26", 32000, 1600000, 20000, 0, 100)
Set-Row 8 @("NH", "2024-02-26", "케이엔알시스템", "DB, NH", "DB, NH", "2024-02-29", "2024-03-07", 14202, 2104000, 13500, 0, 50)

# --- Insert the new 하나32호스팩 row at position 17 ---
# This pushes the existing rows 17 (한국/삼현) and 18 (한화/이에이트) down to
# 18 and 19 respectively, with their original values intact.
$ws.Rows("17:17").Insert()

$ws.Range("B17").NumberFormat = "@"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
Set-Row 17 @("하나", "2024-03-18", "하나32호스팩", "하나", "하나", "2024-03-21", "2024-03-27", 6000, 3000000, 2000, 0, 100)

# Restore default/general number format on the date columns we touched.
$ws.Range("B2:B19").Style = "Normal"
$ws.Range("F2:F19").Style = "Normal"
$ws.Range("G2:G19").Style = "Normal"
